$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 572.375
$ws.Range("I6").Value = 847.8
$ws.Range("K6").Value = 2543.4
$ws.Range("M6").Value = -2431.4
$ws.Range("H11").Value = 20.083334
$ws.Range("I11").Value = 20.083334
$ws.Range("K11").Value = 20.083334
$ws.Range("M11").Value = 119.916666
$ws.Range("H32").Value = 5292
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H38").Value = 10972.111
$ws.Range("I38").Value = 249.8
$ws.Range("J38").Value = 24375
$ws.Range("K38").Value = 749.4000000000001
$ws.Range("L38").Value = 73125
$ws.Range("M38").Value = -377.4000000000001
$ws.Range("N38").Value = -73869
$ws.Range("H39").Value = 1348.4286
$ws.Range("I39").Value = 1208
$ws.Range("J39").Value = 1699.5
$ws.Range("K39").Value = 3624
$ws.Range("L39").Value = 5098.5
$ws.Range("M39").Value = -3328
$ws.Range("N39").Value = -5690.5
$ws.Range("H40").Value = 6146.4443
$ws.Range("I40").Value = 7930
$ws.Range("J40").Value = 2579.3333
$ws.Range("K40").Value = 7930
$ws.Range("L40").Value = 2579.3333
$ws.Range("M40").Value = -7755
$ws.Range("N40").Value = -2929.3333
$ws.Range("H41").Value = 2087
$ws.Range("J41").Value = 2217.111
$ws.Range("L41").Value = 2217.111
$ws.Range("N41").Value = -3097.111
$ws.Range("H42").Value = 385.16666
$ws.Range("I42").Value = 430.5
$ws.Range("J42").Value = 294.5
$ws.Range("K42").Value = 1291.5
$ws.Range("L42").Value = 883.5
$ws.Range("M42").Value = -1061.5
$ws.Range("N42").Value = -1343.5
$ws.Range("H43").Value = 3498.5
$ws.Range("I43").Value = 2999
$ws.Range("K43").Value = 2999
$ws.Range("M43").Value = -2930
$ws.Range("H51").Value = 7575.1113
$ws.Range("I51").Value = 6298.75
$ws.Range("J51").Value = 8596.200000000001
$ws.Range("K51").Value = 6298.75
$ws.Range("L51").Value = 8596.200000000001
$ws.Range("M51").Value = -5814.75
$ws.Range("N51").Value = -9564.200000000001
$ws.Range("H52").Value = 2500
$ws.Range("I52").Value = 2500
$ws.Range("K52").Value = 7500
$ws.Range("M52").Value = -7340
$ws.Range("H57").Value = 157500
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 157500
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 472500
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -473498
$ws.Range("H58").Value = 1508.1428
$ws.Range("I58").Value = 150
$ws.Range("J58").Value = 1734.5
$ws.Range("K58").Value = 450
$ws.Range("L58").Value = 5203.5
$ws.Range("M58").Value = -300
$ws.Range("N58").Value = -5503.5
$ws.Range("H64").Value = 5858.8184
$ws.Range("I64").Value = 5682.75
$ws.Range("K64").Value = 5682.75
$ws.Range("M64").Value = -5434.75
$ws.Range("H67").Value = 5858.8184
$ws.Range("I67").Value = 5682.75
$ws.Range("K67").Value = 5682.75
$ws.Range("M67").Value = -4824.75
$ws.Range("H70").Value = 11855.5
$ws.Range("I70").Value = 4733
$ws.Range("J70").Value = 13798
$ws.Range("K70").Value = 14199
$ws.Range("L70").Value = 41394
$ws.Range("M70").Value = -13929
$ws.Range("N70").Value = -41934
$ws.Range("H73").Value = 11855.5
$ws.Range("I73").Value = 4733
$ws.Range("J73").Value = 13798
$ws.Range("K73").Value = 14199
$ws.Range("L73").Value = 41394
$ws.Range("M73").Value = -13263
$ws.Range("N73").Value = -43266
$ws.Range("H74").Value = 4859.7856
$ws.Range("I74").Value = 4821.5454
$ws.Range("K74").Value = 4821.5454
$ws.Range("M74").Value = -3885.5454
$ws.Range("H77").Value = 4859.7856
$ws.Range("I77").Value = 4821.5454
$ws.Range("K77").Value = 24107.727
$ws.Range("M77").Value = -19427.727
$ws.Range("H86").Value = 57700.61
$ws.Range("I86").Value = 145255.72
$ws.Range("K86").Value = 145255.72
$ws.Range("M86").Value = -144132.72
$ws.Range("H89").Value = 57700.61
$ws.Range("I89").Value = 145255.72
$ws.Range("K89").Value = 726278.6
$ws.Range("M89").Value = -720662.6
$ws.Range("H107").Value = 4031.3809
$ws.Range("I107").Value = 3405.6428
$ws.Range("J107").Value = 5282.857
$ws.Range("K107").Value = 3405.6428
$ws.Range("L107").Value = 5282.857
$ws.Range("M107").Value = -1485.6428
$ws.Range("N107").Value = -9122.857
$ws.Range("H112").Value = 3321078
$ws.Range("J112").Value = 3874257
$ws.Range("L112").Value = 11622771
$ws.Range("N112").Value = -11624987
$ws.Range("H132").Value = 1625.258
$ws.Range("I132").Value = 1580.0714
$ws.Range("K132").Value = 4740.2142
$ws.Range("M132").Value = -2210.2142
$ws.Range("H138").Value = 3471.8687
$ws.Range("I138").Value = 1987.4
$ws.Range("J138").Value = 3973.3784
$ws.Range("K138").Value = 5962.200000000001
$ws.Range("L138").Value = 11920.1352
$ws.Range("M138").Value = -822.2000000000007
$ws.Range("N138").Value = -22200.1352
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14360.706
$ws.Range("I32").Value = 14498
$ws.Range("J32").Value = 10997
$ws.Range("K32").Value = 14498
$ws.Range("L32").Value = 10997
$ws.Range("M32").Value = -14211
$ws.Range("N32").Value = -11571
$ws.Range("H45").Value = 1930
$ws.Range("I45").Value = 1203.8334
$ws.Range("J45").Value = 2474.625
$ws.Range("K45").Value = 1203.8334
$ws.Range("L45").Value = 2474.625
$ws.Range("M45").Value = -826.8334
$ws.Range("N45").Value = -3228.625
$ws.Range("H74").Value = 73297.32000000001
$ws.Range("I74").Value = 73297.32000000001
$ws.Range("K74").Value = 73297.32000000001
$ws.Range("M74").Value = -72423.32000000001
$ws.Range("H77").Value = 73297.32000000001
$ws.Range("I77").Value = 73297.32000000001
$ws.Range("K77").Value = 366486.6
$ws.Range("M77").Value = -362118.6
$ws.Range("H88").Value = 1034.5151
$ws.Range("J88").Value = 951.8570999999999
$ws.Range("L88").Value = 951.8570999999999
$ws.Range("N88").Value = -1763.8571
$ws.Range("H91").Value = 1034.5151
$ws.Range("J91").Value = 951.8570999999999
$ws.Range("L91").Value = 951.8570999999999
$ws.Range("N91").Value = -3759.8571
$ws.Range("H101").Value = 356667
$ws.Range("J101").Value = 356667
$ws.Range("L101").Value = 356667
$ws.Range("N101").Value = -363157
$ws.Range("H122").Value = 9138.879999999999
$ws.Range("I122").Value = 6603.0557
$ws.Range("K122").Value = 19809.1671
$ws.Range("M122").Value = -17359.1671
$ws.Range("H132").Value = 7391.4136
$ws.Range("I132").Value = 5518.769
$ws.Range("K132").Value = 16556.307
$ws.Range("M132").Value = -14026.307
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2321.1892
$ws.Range("I20").Value = 1700.4642
$ws.Range("J20").Value = 4252.3335
$ws.Range("K20").Value = 1700.4642
$ws.Range("L20").Value = 4252.3335
$ws.Range("M20").Value = -1453.4642
$ws.Range("N20").Value = -4746.3335
$ws.Range("H86").Value = 1567.7727
$ws.Range("I86").Value = 1444
$ws.Range("K86").Value = 1444
$ws.Range("M86").Value = -321
$ws.Range("H89").Value = 1567.7727
$ws.Range("I89").Value = 1444
$ws.Range("K89").Value = 7220
$ws.Range("M89").Value = -1604
$ws.Range("H99").Value = 7262.0356
$ws.Range("I99").Value = 8166.85
$ws.Range("J99").Value = 5000
$ws.Range("K99").Value = 8166.85
$ws.Range("L99").Value = 5000
$ws.Range("M99").Value = -6668.85
$ws.Range("N99").Value = -7996
$ws.Range("H105").Value = 1632.4584
$ws.Range("I105").Value = 1892.6428
$ws.Range("J105").Value = 1268.2
$ws.Range("K105").Value = 1892.6428
$ws.Range("L105").Value = 1268.2
$ws.Range("M105").Value = -145.6428000000001
$ws.Range("N105").Value = -4762.2
$ws.Range("H134").Value = 1666.841
$ws.Range("I134").Value = 1239.5135
$ws.Range("K134").Value = 3718.5405
$ws.Range("M134").Value = -1183.5405
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1301.3478
$ws.Range("I16").Value = 1184.7894
$ws.Range("J16").Value = 1855
$ws.Range("K16").Value = 1184.7894
$ws.Range("L16").Value = 1855
$ws.Range("M16").Value = -897.7893999999999
$ws.Range("N16").Value = -2429
$ws.Range("H58").Value = 21167.625
$ws.Range("I58").Value = 5192.8
$ws.Range("J58").Value = 47792.332
$ws.Range("K58").Value = 5192.8
$ws.Range("L58").Value = 47792.332
$ws.Range("M58").Value = -4989.8
$ws.Range("N58").Value = -48198.332
$ws.Range("H99").Value = 2003161.2
$ws.Range("J99").Value = 3947
$ws.Range("L99").Value = 3947
$ws.Range("N99").Value = -6943
$ws.Range("H113").Value = 1301.3478
$ws.Range("I113").Value = 1184.7894
$ws.Range("J113").Value = 1855
$ws.Range("K113").Value = 1184.7894
$ws.Range("L113").Value = 1855
$ws.Range("M113").Value = 985.2106000000001
$ws.Range("N113").Value = -6195
$ws.Range("H126").Value = 2003161.2
$ws.Range("J126").Value = 3947
$ws.Range("L126").Value = 11841
$ws.Range("N126").Value = -16781
$ws.Range("H132").Value = 4868.1665
$ws.Range("I132").Value = 2906.182
$ws.Range("K132").Value = 8718.545999999998
$ws.Range("M132").Value = -6188.545999999998
$ws.Range("H134").Value = 6994.8096
$ws.Range("I134").Value = 7725.8335
$ws.Range("K134").Value = 23177.5005
$ws.Range("M134").Value = -20642.5005
$ws.Range("H136").Value = 21167.625
$ws.Range("I136").Value = 5192.8
$ws.Range("J136").Value = 47792.332
$ws.Range("K136").Value = 15578.4
$ws.Range("L136").Value = 143376.996
$ws.Range("M136").Value = -13028.4
$ws.Range("N136").Value = -148476.996
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 297
$ws.Range("I9").Value = 297.25
$ws.Range("K9").Value = 891.75
$ws.Range("M9").Value = -667.75
$ws.Range("H41").Value = 2000
$ws.Range("J41").Value = 2000
$ws.Range("L41").Value = 6000
$ws.Range("N41").Value = -6676
$ws.Range("H43").Value = 2003
$ws.Range("J43").Value = 2003
$ws.Range("L43").Value = 6009
$ws.Range("N43").Value = -6237
$ws.Range("H117").Value = 199.33333
$ws.Range("I117").Value = 199.33333
$ws.Range("K117").Value = 597.99999
$ws.Range("M117").Value = 2844.00001
$ws.Range("H121").Value = 2679.8
$ws.Range("I121").Value = 1439.75
$ws.Range("K121").Value = 4319.25
$ws.Range("M121").Value = -3009.25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 14999.25
$ws.Range("J20").Value = 14999.25
$ws.Range("L20").Value = 14999.25
$ws.Range("N20").Value = -15489.25
$ws.Range("H39").Value = 42500
$ws.Range("J39").Value = 42500
$ws.Range("L39").Value = 42500
$ws.Range("N39").Value = -43564
$ws.Range("H102").Value = 2471.1667
$ws.Range("I102").Value = 1085.1786
$ws.Range("K102").Value = 1085.1786
$ws.Range("M102").Value = 536.8214
$ws.Range("H132").Value = 43597.168
$ws.Range("I132").Value = 56912.4
$ws.Range("J132").Value = 16966.7
$ws.Range("K132").Value = 170737.2
$ws.Range("L132").Value = 50900.10000000001
$ws.Range("M132").Value = -168207.2
$ws.Range("N132").Value = -55960.10000000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 520.6818
$ws.Range("I16").Value = 520.6818
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 520.6818
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -350.6818
$ws.Range("N16").ClearContents()
$ws.Range("H55").Value = 9690.941000000001
$ws.Range("I55").Value = 748.5
$ws.Range("J55").Value = 31152.8
$ws.Range("K55").Value = 748.5
$ws.Range("L55").Value = 31152.8
$ws.Range("M55").Value = -575.5
$ws.Range("N55").Value = -31498.8
$ws.Range("H61").Value = 1197.8
$ws.Range("I61").Value = 1197.8
$ws.Range("K61").Value = 1197.8
$ws.Range("M61").Value = -995.8
$ws.Range("H113").Value = 1197.8
$ws.Range("I113").Value = 1197.8
$ws.Range("K113").Value = 1197.8
$ws.Range("M113").Value = 972.2
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H123").Value = 19999
$ws.Range("J123").Value = 19999
$ws.Range("L123").Value = 19999
$ws.Range("N123").Value = -29799
$ws.Range("H132").Value = 7158.08
$ws.Range("I132").Value = 3448.7896
$ws.Range("K132").Value = 10346.3688
$ws.Range("M132").Value = -7816.3688
$ws.Range("H136").Value = 6468.778
$ws.Range("I136").Value = 6839.875
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 20519.625
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -17969.625
$ws.Range("N136").Value = -15600
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H37").Value = 79831.336
$ws.Range("H52").Value = 37666.332
$ws.Range("I52").Value = 44500
$ws.Range("K52").Value = 44500
$ws.Range("M52").Value = -44274
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H107").Value = 1198.8
$ws.Range("I107").Value = 1198.5
$ws.Range("J107").Value = 1199
$ws.Range("K107").Value = 3595.5
$ws.Range("L107").Value = 3597
$ws.Range("M107").Value = -1675.5
$ws.Range("N107").Value = -7437
$ws.Range("H113").Value = 1755.6154
$ws.Range("I113").Value = 1747
$ws.Range("J113").Value = 1784.3334
$ws.Range("K113").Value = 5241
$ws.Range("L113").Value = 5353.0002
$ws.Range("M113").Value = -3071
$ws.Range("N113").Value = -9693.0002
$ws.Range("H122").Value = 2750
$ws.Range("I122").Value = 2750
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8250
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5800
$ws.Range("N122").ClearContents()
$ws.Range("H136").Value = 178214.94
$ws.Range("I136").Value = 239990.73
$ws.Range("K136").Value = 719972.1900000001
$ws.Range("M136").Value = -717422.1900000001
